$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "1202"
$t.Cell(5, 1).Range.Text = "0.00001"
$t.Cell(6, 1).Range.Text = "0.00058"
$t.Cell(8, 1).Range.Text = "0.00004"
$t.Cell(9, 1).Range.Text = "0.00028"
$t.Cell(10, 1).Range.Text = "0.00036"
$t.Cell(11, 1).Range.Text = "0.00039"
$t.Cell(12, 1).Range.Text = "0.23149"

$t.Cell(44, 1).Range.Text = "99.93"
$t.Cell(45, 1).Range.Text = "0.23"
$t.Cell(46, 1).Range.Text = "340"
